# Add a new trailing row (row 80) to each of the four sensor-log sheets,
# mirroring the existing row 79 in each sheet but with the hour advanced
# by one (14:xx:xx -> 15:xx:xx) as captured in the target OOXML diff.

function Set-RowData($ws, $row, $a, $b, $c, $d, $e, $f, $g, $h, $i) {
    $ws.Range("A$row").Value = $a
    $ws.Range("B$row").Value = $b
    $ws.Range("C$row").Value = $c
    $ws.Range("D$row").Value = $d
    $ws.Range("E$row").Value = $e
    $ws.Range("F$row").Value = $f

    # G holds a 24-digit numeric string that must stay text (it overflows
    # float precision); force text entry, then drop back to the sheet's
    # default "Normal" style so no stray numFmt is left on the cell.
    $ws.Range("G$row").NumberFormat = "@"
    $ws.Range("G$row").Value = $g
    $ws.Range("G$row").Style = "Normal"

    $ws.Range("H$row").Value = $h
    $ws.Range("I$row").Value = $i
}

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ROW35-FE-LIFTER")
Set-RowData $ws1 80 "2025-03-07 15:42:06" "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c," "0x01,0x90," "0x d" 400 "568631262647113770877196" 400 13

$ws2 = $wb.Worksheets.Item("ROW35-MID-LIFTER")
Set-RowData $ws2 80 "2025-03-07 15:29:35" "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c," "0x01,0x90," "0x e" 400 "568631262647113770942732" 400 14

$ws3 = $wb.Worksheets.Item("ROW02-FE-LIFTER")
Set-RowData $ws3 80 "2025-03-07 15:51:45" "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c," "0x01,0x90," "0xff" 400 "568631262647113769959692" 400 255

$ws4 = $wb.Worksheets.Item("ROW02-MID-LIFTER")
Set-RowData $ws4 80 "2025-03-07 15:41:15" "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c," "0x01,0x90," "0x 3" 400 "568631262647113769959692" 400 3
